$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Load a CSV-style map into rows 90-92 (column A keeps the default
# style, columns B:CL get the "0.0" numeric format applied with an
# explicit font touch, matching the new cellXfs entry) ---

# Column A (index 1) values for rows 90:92 stay on the original style,
# just make sure every row in the block has a value (row 91 is brand new,
# so its style has to be (re)stamped explicitly to line up with A90/A92).
$ws.Range("A90:A92").Value = 0
$ws.Range("A90:A92").NumberFormat = "0.0"

# Columns B:CL (2:90) for rows 90:92 get the numeric format re-applied
# (this is what stamps the extra cellXfs entry used by the new rows).
$dataRange = $ws.Range("B90:CL92")
$dataRange.Value = 0
$dataRange.NumberFormat = "0.0"
$dataRange.Font.Name = "Calibri"

# --- Conditional formatting now covers the extended block (rows 3-92) ---
$cf = $ws.Range("A3:CL90").FormatConditions
$cf.Item(1).ModifyAppliesToRange($ws.Range("A3:CL92"))

# --- Update the view: scrolled to A68, active selection on V92 ---
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 68
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("V92").Select() | Out-Null
